$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.432.17"
$ws.Range("E2").Value = "  -3.00%  "

$ws.Range("D3").Value = "'1.802.36"
$ws.Range("E3").Value = "  -2.70%  "

$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.55%  "

$ws.Range("D6").Value = "'307.97"
$ws.Range("E6").Value = "  -1.88%  "

$ws.Range("D7").Value = "'0.4528"
$ws.Range("E7").Value = "  -1.57%  "

$ws.Range("E8").Value = "  -1.75%  "

$ws.Range("D9").Value = "'0.07091"
$ws.Range("E9").Value = "  -2.75%  "

$ws.Range("D10").Value = "'0.8724"
$ws.Range("E10").Value = "  -1.62%  "

$ws.Range("D11").Value = "'0.07755"
$ws.Range("E11").Value = "  -0.93%  "

$ws.Range("D12").Value = "'19.28"
$ws.Range("E12").Value = "  -4.04%  "

$ws.Range("D13").Value = "'1.812.85"
$ws.Range("E13").Value = "  -0.98%  "

$ws.Range("D14").Value = "'5.253"
$ws.Range("E14").Value = "  -2.52%  "

$ws.Range("D15").Value = "'6.327"
$ws.Range("E15").Value = "  -3.08%  "

$ws.Range("D16").Value = "'85.86"
$ws.Range("E16").Value = "  -6.05%  "

$ws.Range("D17").Value = "'1.009"
$ws.Range("E17").Value = "  +0.50%  "

$ws.Range("D18").Value = "'0.000008557"
$ws.Range("E18").Value = "  -4.20%  "

$ws.Range("D19").Value = "'1.007"
$ws.Range("E19").Value = "  +0.55%  "

$ws.Range("D20").Value = "'26.494.53"
$ws.Range("E20").Value = "  -2.88%  "

$ws.Range("D21").Value = "'14.20"
$ws.Range("E21").Value = "  -3.75%  "

$ws.Range("D22").Value = "'4.955"
$ws.Range("E22").Value = "  -3.04%  "

$ws.Range("D23").Value = "'10.36"
$ws.Range("E23").Value = "  -1.63%  "

$ws.Range("D24").Value = "'1.968"
$ws.Range("E24").Value = "  +2.71%  "

$ws.Range("D25").Value = "'150.58"
$ws.Range("E25").Value = "  -0.90%  "

$ws.Range("D26").Value = "'17.86"
$ws.Range("E26").Value = "  -3.12%  "

$ws.Range("D27").Value = "'1.991"
$ws.Range("E27").Value = "  -3.23%  "

$ws.Range("D28").Value = "'112.83"
$ws.Range("E28").Value = "  -2.59%  "

$ws.Range("D29").Value = "'4.845"
$ws.Range("E29").Value = "  -4.28%  "

$ws.Range("D30").Value = "'0.08647"
$ws.Range("E30").Value = "  -1.99%  "

$ws.Range("D31").Value = "'3.033"
$ws.Range("E31").Value = "  -1.31%  "

$ws.Range("D32").Value = "'0.7272"
$ws.Range("E32").Value = "  -5.91%  "

$ws.Range("D33").Value = "'4.436"
$ws.Range("E33").Value = "  -1.31%  "

$ws.Range("D34").Value = "'1.108"
$ws.Range("E34").Value = "  -5.18%  "

$ws.Range("D35").Value = "'1.005"
$ws.Range("E35").Value = "  +0.49%  "

$ws.Range("D36").Value = "'2.528"
$ws.Range("E36").Value = "  -7.87%  "

$ws.Range("E37").Value = "  -1.10%  "

$ws.Range("D38").Value = "'0.01921"
$ws.Range("E38").Value = "  -1.53%  "

$ws.Range("D39").Value = "'0.05067"
$ws.Range("E39").Value = "  -3.67%  "

$ws.Range("E40").Value = "  -2.81%  "

$ws.Range("D41").Value = "'6.944"
$ws.Range("E41").Value = "  -1.63%  "

$ws.Range("D42").Value = "'0.4946"
$ws.Range("E42").Value = "  -3.40%  "

$ws.Range("D43").Value = "'0.1563"
$ws.Range("E43").Value = "  -4.35%  "

$ws.Range("D44").Value = "'8.094"
$ws.Range("E44").Value = "  -3.56%  "

$ws.Range("E45").Value = "  +0.56%  "

$ws.Range("D46").Value = "'0.4599"
$ws.Range("E46").Value = "  -4.10%  "

$ws.Range("D47").Value = "'101.52"
$ws.Range("E47").Value = "  -0.80%  "

$ws.Range("D48").Value = "'9.914"
$ws.Range("E48").Value = "  -4.13%  "

$ws.Range("E49").Value = "  -3.88%  "

$ws.Range("D50").Value = "'0.05984"
$ws.Range("E50").Value = "  -3.72%  "

$ws.Range("D51").Value = "'63.53"
$ws.Range("E51").Value = "  -3.32%  "
